# Book1.xlsx / Sheet2 ("se" class) — append 4 more tracked-instance rows
# below the existing "goo" rows, mirroring the same layout
# (Symbol | Date | Price), then move the selection to the next empty row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 10; Date = 41244 },
    @{ Row = 11; Date = 41245 },
    @{ Row = 12; Date = 41246 },
    @{ Row = 13; Date = 41247 }
)

foreach ($r in $rows) {
    $ws.Range("A$($r.Row)").Value = "se"
    $ws.Range("B$($r.Row)").Value = $r.Date
    $ws.Range("B$($r.Row)").NumberFormat = "m/d/yy"
    $ws.Range("C$($r.Row)").Value = 40
}

$ws.Range("A14").Select()
